$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increase hole diameter from 0.75 mm to 0.90 mm
$ws.Range("B21").Value = "0.90 mm"

# Update selection to reflect where the user ended up after the edit
$ws.Range("D23").Select()
